$wb = $excel.ActiveWorkbook

# ALC sheet changes (@@ -1477,25 +1477,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 652.0421
$ws.Range("J17").Value = 555.55
$ws.Range("L17").Value = 1666.65
$ws.Range("N17").Value = -2002.65

# ALC sheet changes (@@ -5820,25 +5820,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 700
$ws.Range("I103").Value = 700
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 2100
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -1514
$ws.Range("N103").ClearContents()

# ALC sheet changes (@@ -6267,25 +6264,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1583.4182
$ws.Range("I112").Value = 1124
$ws.Range("J112").Value = 1600.7548
$ws.Range("K112").Value = 3372
$ws.Range("L112").Value = 4802.2644
$ws.Range("M112").Value = -2264
$ws.Range("N112").Value = -7018.2644

# ALC sheet changes (@@ -6466,22 +6463,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 536290.0600000001
$ws.Range("I116").Value = 2503701.2
$ws.Range("K116").Value = 2503701.2
$ws.Range("M116").Value = -2500259.2

# ALC sheet changes (@@ -6567,22 +6564,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 697.7778
$ws.Range("I118").Value = 415.55554
$ws.Range("K118").Value = 1246.66662
$ws.Range("M118").Value = 410.33338

# ALC sheet changes (@@ -7265,25 +7262,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 23492990
$ws.Range("I132").Value = 29416452
$ws.Range("J132").Value = 1115473.5
$ws.Range("K132").Value = 88249356
$ws.Range("L132").Value = 3346420.5
$ws.Range("M132").Value = -88246826
$ws.Range("N132").Value = -3351480.5

# ALC sheet changes (@@ -7516,25 +7513,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2618.432
$ws.Range("I137").Value = 1269.5
$ws.Range("J137").Value = 3967.3635
$ws.Range("K137").Value = 3808.5
$ws.Range("L137").Value = 11902.0905
$ws.Range("M137").Value = -1258.5
$ws.Range("N137").Value = -17002.0905

# ALC sheet changes (@@ -7568,25 +7565,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3161.08
$ws.Range("I138").Value = 636.9524
$ws.Range("J138").Value = 4988.8965
$ws.Range("K138").Value = 1910.8572
$ws.Range("L138").Value = 14966.6895
$ws.Range("M138").Value = 3229.1428
$ws.Range("N138").Value = -25246.6895

# ALC sheet changes (@@ -7721,25 +7718,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3006.7097
$ws.Range("I141").Value = 2804.453
$ws.Range("J141").Value = 4197.778
$ws.Range("K141").Value = 8413.359
$ws.Range("L141").Value = 12593.334
$ws.Range("M141").Value = -3233.359
$ws.Range("N141").Value = -22953.334

# ARM sheet changes (@@ -9355,25 +9352,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3500.0618
$ws.Range("I32").Value = 3207.139
$ws.Range("J32").Value = 5843.4443
$ws.Range("K32").Value = 3207.139
$ws.Range("L32").Value = 5843.4443
$ws.Range("M32").Value = -2920.139
$ws.Range("N32").Value = -6417.4443

# ARM sheet changes (@@ -10794,25 +10791,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 819.6
$ws.Range("I61").Value = 586.5185
$ws.Range("J61").Value = 1606.25
$ws.Range("K61").Value = 586.5185
$ws.Range("L61").Value = 1606.25
$ws.Range("M61").Value = -374.5185
$ws.Range("N61").Value = -2030.25

# ARM sheet changes (@@ -11422,25 +11419,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2684.6274
$ws.Range("I74").Value = 2504.6223
$ws.Range("J74").Value = 4034.6667
$ws.Range("K74").Value = 2504.6223
$ws.Range("L74").Value = 4034.6667
$ws.Range("M74").Value = -1630.6223
$ws.Range("N74").Value = -5782.6667

# ARM sheet changes (@@ -11572,25 +11569,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2684.6274
$ws.Range("I77").Value = 2504.6223
$ws.Range("J77").Value = 4034.6667
$ws.Range("K77").Value = 12523.1115
$ws.Range("L77").Value = 20173.3335
$ws.Range("M77").Value = -8155.111499999999
$ws.Range("N77").Value = -28909.3335

# ARM sheet changes (@@ -12543,25 +12540,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 654.3333
$ws.Range("I97").Value = 703.86664
$ws.Range("J97").Value = 406.66666
$ws.Range("K97").Value = 703.86664
$ws.Range("L97").Value = 406.66666
$ws.Range("M97").Value = -207.86664
$ws.Range("N97").Value = -1398.66666

# ARM sheet changes (@@ -14237,25 +14234,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2616.2856
$ws.Range("I132").Value = 1617.3928
$ws.Range("J132").Value = 4614.0713
$ws.Range("K132").Value = 4852.178400000001
$ws.Range("L132").Value = 13842.2139
$ws.Range("M132").Value = -2322.178400000001
$ws.Range("N132").Value = -18902.2139

# ARM sheet changes (@@ -14436,25 +14433,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 819.6
$ws.Range("I136").Value = 586.5185
$ws.Range("J136").Value = 1606.25
$ws.Range("K136").Value = 1759.5555
$ws.Range("L136").Value = 4818.75
$ws.Range("M136").Value = 790.4445000000001
$ws.Range("N136").Value = -9918.75

# BSM sheet changes (@@ -21250,25 +21247,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1529.7727
$ws.Range("I134").Value = 907.6667
$ws.Range("J134").Value = 3396.0908
$ws.Range("K134").Value = 2723.0001
$ws.Range("L134").Value = 10188.2724
$ws.Range("M134").Value = -188.0001000000002
$ws.Range("N134").Value = -15258.2724

# CRP sheet changes (@@ -23178,25 +23175,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8199373.5
$ws.Range("I31").Value = 1458.6216
$ws.Range("J31").Value = 20837826
$ws.Range("K31").Value = 1458.6216
$ws.Range("L31").Value = 20837826
$ws.Range("M31").Value = -1163.6216
$ws.Range("N31").Value = -20838416

# CRP sheet changes (@@ -23331,25 +23328,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8199373.5
$ws.Range("I34").Value = 1458.6216
$ws.Range("J34").Value = 20837826
$ws.Range("K34").Value = 1458.6216
$ws.Range("L34").Value = 20837826
$ws.Range("M34").Value = -1256.6216
$ws.Range("N34").Value = -20838230

# CRP sheet changes (@@ -24504,22 +24501,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1753.663
$ws.Range("I58").Value = 1553.8857
$ws.Range("K58").Value = 1553.8857
$ws.Range("M58").Value = -1350.8857

# CRP sheet changes (@@ -26495,25 +26492,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 9529207
$ws.Range("I99").Value = 15388411
$ws.Range("J99").Value = 8000.5
$ws.Range("K99").Value = 15388411
$ws.Range("L99").Value = 8000.5
$ws.Range("M99").Value = -15386913
$ws.Range("N99").Value = -10996.5

# CRP sheet changes (@@ -27812,25 +27809,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 9529207
$ws.Range("I126").Value = 15388411
$ws.Range("J126").Value = 8000.5
$ws.Range("K126").Value = 46165233
$ws.Range("L126").Value = 24001.5
$ws.Range("M126").Value = -46162763
$ws.Range("N126").Value = -28941.5

# CRP sheet changes (@@ -28109,25 +28106,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2403.158
$ws.Range("I132").Value = 1881.3549
$ws.Range("J132").Value = 4714
$ws.Range("K132").Value = 5644.0647
$ws.Range("L132").Value = 14142
$ws.Range("M132").Value = -3114.0647
$ws.Range("N132").Value = -19202

# CRP sheet changes (@@ -28210,25 +28207,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3137.2205
$ws.Range("I134").Value = 3676.6572
$ws.Range("J134").Value = 2350.5417
$ws.Range("K134").Value = 11029.9716
$ws.Range("L134").Value = 7051.625100000001
$ws.Range("M134").Value = -8494.971600000001
$ws.Range("N134").Value = -12121.6251

# CRP sheet changes (@@ -28311,22 +28308,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1753.663
$ws.Range("I136").Value = 1553.8857
$ws.Range("K136").Value = 4661.6571
$ws.Range("M136").Value = -2111.6571

# CUL sheet changes (@@ -32044,25 +32041,25 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 13472.125
$ws.Range("I68").Value = 735
$ws.Range("J68").Value = 17717.834
$ws.Range("K68").Value = 2205
$ws.Range("L68").Value = 53153.50199999999
$ws.Range("M68").Value = -1394
$ws.Range("N68").Value = -54775.50199999999

# CUL sheet changes (@@ -32200,25 +32197,25 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 13472.125
$ws.Range("I71").Value = 735
$ws.Range("J71").Value = 17717.834
$ws.Range("K71").Value = 6615
$ws.Range("L71").Value = 159460.506
$ws.Range("M71").Value = -2559
$ws.Range("N71").Value = -167572.506

# CUL sheet changes (@@ -35285,25 +35282,25 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1880.1333
$ws.Range("I132").Value = 754.7143
$ws.Range("J132").Value = 2864.875
$ws.Range("K132").Value = 6792.428699999999
$ws.Range("L132").Value = 25783.875
$ws.Range("M132").Value = -4262.428699999999
$ws.Range("N132").Value = -30843.875

# GSM sheet changes (@@ -42236,22 +42233,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4974.643
$ws.Range("I132").Value = 3442
$ws.Range("K132").Value = 10326
$ws.Range("M132").Value = -7796

# LTW sheet changes (@@ -47624,22 +47621,22 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H99").Value = 34950
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 34950
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 34950
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -40940

# WVR sheet changes (@@ -52538,22 +52535,22 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 30075
$ws.Range("J57").Value = 30075
$ws.Range("L57").Value = 30075
$ws.Range("N57").Value = -31583

# WVR sheet changes (@@ -55711,22 +55708,22 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4378.9375
$ws.Range("I122").Value = 2896.4443
$ws.Range("K122").Value = 8689.332900000001
$ws.Range("M122").Value = -6239.332900000001

# WVR sheet changes (@@ -56204,25 +56201,25 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6174285
$ws.Range("I132").Value = 697.6774
$ws.Range("J132").Value = 14495207
$ws.Range("K132").Value = 2093.0322
$ws.Range("L132").Value = 43485621
$ws.Range("M132").Value = 436.9677999999999
$ws.Range("N132").Value = -43490681

# WVR sheet changes (@@ -56403,25 +56400,25 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1714.6165
$ws.Range("I136").Value = 617.2341
$ws.Range("J136").Value = 3698.3462
$ws.Range("K136").Value = 1851.7023
$ws.Range("L136").Value = 11095.0386
$ws.Range("M136").Value = 698.2977000000001
$ws.Range("N136").Value = -16195.0386

